$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62
$ws.Range("A62").Value = "aayayy"
$ws.Range("B62").Value = 55
$ws.Range("C62").Value = "asdfasdf"
$ws.Range("D62").Value = "asdfasdf"

# Row 63
$ws.Range("A63").Value = "last"
$ws.Range("B63").Value = 2
$ws.Range("C63").Value = "'22"
$ws.Range("C63").ClearFormats()
$ws.Range("D63").Value = "'22"
$ws.Range("D63").ClearFormats()

# Rows 64 and 65 remain blank but are touched so they stay part of the
# worksheet's used range (matching the new dimension A1:D65).
$ws.Rows.Item(64).Hidden = $true
$ws.Rows.Item(64).Hidden = $false
$ws.Rows.Item(65).Hidden = $true
$ws.Rows.Item(65).Hidden = $false
